$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 4: D4 becomes a number (1) instead of text "1,0"
$ws2.Range("D4").Value = 1

# Row 5: C5 becomes empty (was "Wert2"), D5 becomes a number (1.1) instead of text "1,1"
$ws2.Range("C5").Value = ""
$ws2.Range("D5").Value = 1.1

# Row 6: C6 becomes empty (was "Wert3"), D6 becomes a number (1.2) instead of text "1,2"
$ws2.Range("C6").Value = ""
$ws2.Range("D6").Value = 1.2

# Row 7: D7 becomes empty (was "1,3")
$ws2.Range("D7").Value = ""

# Row 8: D8 becomes a number (1.4) instead of text "1,4"
$ws2.Range("D8").Value = 1.4
